$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.912.12"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "1.630.21"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.88"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("E6").Value = "  -1.36%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.24"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0883"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").Value = "1.860.87"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "1.631.24"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("D17").Value = "27.913.74"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.33"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.35"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.00"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.87"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.92"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("E29").Value = "  -1.13%  "
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0481"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").Value = "1.415.37"
$ws.Range("E33").Value = "  +0.89%  "
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("E35").Value = "  +2.55%  "
$ws.Range("E36").Value = "  -1.71%  "
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0170"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.18%  "
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("E40").Value = "  -1.83%  "
$ws.Range("E41").Value = "  -2.06%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.83"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.48%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.77"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("E44").Value = "  -0.87%  "
$ws.Range("D45").Value = "1.770.70"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("E46").Value = "  -3.80%  "
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("E48").Value = "  -2.46%  "
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.59"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.70%  "
